$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.355.82'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.882.48'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.92'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9993'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4832'
$ws.Range('E7').Value = '  -1.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2889'
$ws.Range('E8').Value = '  -2.93%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06595'
$ws.Range('E9').Value = '  -2.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.876.30'
$ws.Range('E10').Value = '  -2.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.92'
$ws.Range('E11').Value = '  -1.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07387'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.172'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.40'
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6605'
$ws.Range('E15').Value = '  -1.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.306.89'
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.49'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007742'
$ws.Range('E18').Value = '  -2.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9989'
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.432'
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.135.10'
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '197.92'
$ws.Range('E23').Value = '  -1.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.171'
$ws.Range('E24').Value = '  -2.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.410'
$ws.Range('E25').Value = '  -2.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.02'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.23'
$ws.Range('E27').Value = '  -2.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.932'
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.432'
$ws.Range('E29').Value = '  -3.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.257'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09134'
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05044'
$ws.Range('E33').Value = '  -4.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7382'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.145'
$ws.Range('E35').Value = '  +2.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.703'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01870'
$ws.Range('E37').Value = '  +1.88%  '
$ws.Range('E38').Value = '  -3.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9143'
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.078'
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.966'
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.67'
$ws.Range('E42').Value = '  -0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4326'
$ws.Range('E43').Value = '  -3.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9995'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.592'
$ws.Range('E45').Value = '  -0.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1344'
$ws.Range('E46').Value = '  -3.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.572'
$ws.Range('E47').Value = '  +9.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '65.20'
$ws.Range('E48').Value = '  -13.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.907'
$ws.Range('E49').Value = '  -2.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.22'
$ws.Range('E50').Value = '  -3.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05684'
$ws.Range('E51').Value = '  -3.26%  '
